# Update the two-digit multiplication problems in the table to the new
# set of operands, per the commit "Update master to output generated at
# c8c62b6".

$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "22×93="; New = "47×88=" },
    @{ Old = "69×76="; New = "72×40=" },
    @{ Old = "65×64="; New = "86×54=" },
    @{ Old = "58×29="; New = "99×89=" },
    @{ Old = "74×55="; New = "68×63=" },
    @{ Old = "66×89="; New = "66×59=" },
    @{ Old = "37×40="; New = "18×40=" },
    @{ Old = "90×72="; New = "77×32=" },
    @{ Old = "27×58="; New = "81×19=" },
    @{ Old = "78×78="; New = "59×57=" },
    @{ Old = "61×97="; New = "24×16=" },
    @{ Old = "88×16="; New = "49×18=" },
    @{ Old = "76×11="; New = "86×24=" },
    @{ Old = "72×35="; New = "87×11=" },
    @{ Old = "88×40="; New = "89×96=" },
    @{ Old = "32×30="; New = "21×91=" },
    @{ Old = "94×84="; New = "67×41=" },
    @{ Old = "22×51="; New = "45×88=" },
    @{ Old = "62×37="; New = "46×32=" },
    @{ Old = "24×81="; New = "35×74=" },
    @{ Old = "72×18="; New = "99×88=" },
    @{ Old = "12×17="; New = "39×62=" },
    @{ Old = "59×79="; New = "14×62=" },
    @{ Old = "75×47="; New = "83×58=" },
    @{ Old = "50×16="; New = "95×65=" }
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)
}
